$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "관리자는 원하는 대여소 항목을 선택하면 등록시 입력한 상세 내용을 볼 수 있다."

$ws.Range("B7").Select()
